$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 9: Inscritos 19 -> 20
$ws.Range("E9").Value = 20

# Row 15: Inscritos 135 -> 136, Pagos 66 -> 67, Inscrições homologadas 66 -> 67
$ws.Range("E15").Value = 136
$ws.Range("F15").Value = 67
$ws.Range("H15").Value = 67

# Row 17: Pagos 37 -> 38, Inscrições homologadas 37 -> 38
$ws.Range("F17").Value = 38
$ws.Range("H17").Value = 38

# Row 33: Inscritos 28 -> 29
$ws.Range("E33").Value = 29

# Row 37: Inscritos 37 -> 39
$ws.Range("E37").Value = 39

# Row 62: Inscritos 33 -> 34
$ws.Range("E62").Value = 34

# Row 68: Inscritos 10 -> 11
$ws.Range("E68").Value = 11

# Row 73: Inscritos 20 -> 21
$ws.Range("E73").Value = 21

$wb.Save()
